$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ A = "Department of Education"; B = "Jim Brown"; E = "United States Senate, Office of Senator Robert P. Casey, Jr. (Retired)"; F = "Volunteer" },
    @{ A = "Department of Justice"; B = "Roosevelt Holmes"; E = "Democratic National Committee"; F = "Transition — PT Fund, Inc." },
    @{ A = "Department of Justice"; B = "Regina “Gina” Kline"; E = "Smartjob LLC"; F = "Volunteer" },
    @{ A = "Department of State"; B = "Uzra Zeya"; E = "Alliance for Peacebuilding"; F = "Volunteer" },
    @{ A = "Executive Office of the President, Management and Administration"; B = "Anthony Bernal"; E = "Biden for President"; F = "Transition — PT Fund, Inc." },
    @{ A = "International Development"; B = "Angelique Crumbly"; E = "United Nations Development Program"; F = "Volunteer" },
    @{ A = "National Security Council"; B = "Monica Maher"; E = "Goldman Sachs & Co."; F = "Volunteer" }
)

$rowIndex = 2
foreach ($entry in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $entry.A
    $ws.Cells.Item($rowIndex, 2).Value = $entry.B
    $ws.Cells.Item($rowIndex, 5).Value = $entry.E
    $ws.Cells.Item($rowIndex, 6).Value = $entry.F
    $rowIndex++
}
